# Add a test of get_BQ to the spreadsheet:
#  - rename Sheet1 -> FOCs
#  - add a new sheet BQ (after FOCs) with a worked example of the BQ (bequest) calc
#  - tweak a couple of FOCs inputs/outputs and mirror column B into column C there

$wb = $excel.ActiveWorkbook

# --- Sheet1 -> FOCs -------------------------------------------------------
$focs = $wb.Worksheets.Item(1)
$focs.Name = "FOCs"

# --- New sheet "BQ", placed right after FOCs ------------------------------
$bq = $wb.Worksheets.Add([Type]::Missing, $focs)
$bq.Name = "BQ"

# ===========================================================================
# FOCs sheet edits
# ===========================================================================

$focs.Range("C7").Value = 0

$focs.Range("C14").Value = 0
$focs.Range("C15").Value = 0.96
$focs.Range("C16").Value = 2

$focs.Range("C18").Formula = "=C10^-C16"
$focs.Range("C19").Formula = "=D10^-C16"
$focs.Range("C21").Formula = "=C18 - C15*(1+D3)*(1-C14)*C19"

$focs.Range("C24").Value = 1
$focs.Range("C25").Value = 0.5
$focs.Range("C26").Value = 0.47
$focs.Range("C27").Value = 1.2

$focs.Range("C29").Formula = "=((C26/C24)*((C5/C24)^(C27-1))*((1-((C5/C24)^C27))^((1-C27)/C27))) * C25"
$focs.Range("C31").Formula = "=C4*C18-C29"

$focs.Range("B29").Select() | Out-Null

# ===========================================================================
# BQ sheet contents
# ===========================================================================

$bq.Range("A2").Value = "b_sp1"
$bq.Range("B2").Value = 0.4
$bq.Range("C2").Value = 0.4
$bq.Range("D2").Value = 0

$bq.Range("B3").Value = 0.3
$bq.Range("C3").Value = 0.5
$bq.Range("D3").Value = 0

$bq.Range("B4").Value = 0.2
$bq.Range("C4").Value = 0.6
$bq.Range("D4").Value = 0

$bq.Range("A7").Value = "omega"
$bq.Range("B7").Value = 0.5
$bq.Range("C7").Value = 0.3
$bq.Range("D7").Value = 0.2
$bq.Range("F7").Formula = "=SUM(B7:D7)"

$bq.Range("B8").Value = 0.5
$bq.Range("C8").Value = 0.3
$bq.Range("D8").Value = 0.2
$bq.Range("F8").Formula = "=SUM(B8:D8)"

$bq.Range("B9").Value = 0.5
$bq.Range("C9").Formula = "=B9*(C8/B8)"
$bq.Range("D9").Formula = "=C9*(D8/C8)"
$bq.Range("F9").Formula = "=SUM(B9:D9)"

$bq.Range("A11").Value = "rho_s"
$bq.Range("B11").Formula = "=1-C7/B7"
$bq.Range("C11").Formula = "=1-(D7/C7)"
$bq.Range("D11").Value = 1

$bq.Range("A13").Value = "r"
$bq.Range("B13").Value = 0.05
$bq.Range("B14").Value = 0.04
$bq.Range("B15").Value = 0.03

$bq.Range("A17").Value = "g_n"
$bq.Range("B17").Value = 0.01
$bq.Range("B18").Value = 0.02
$bq.Range("B19").Value = 0.02

$bq.Range("A21").Value = "BQ"
$bq.Range("B21").Formula = "=((1+B13)/(1+B17))*(B7*B`$11*B2+C2*C7*C`$11+D2*D7*D`$11)"
$bq.Range("B22").Formula = "=((1+B14)/(1+B18))*(B8*B`$11*B3+C3*C8*C`$11+D3*D8*D`$11)"
$bq.Range("B23").Formula = "=((1+B15)/(1+B19))*(B9*B`$11*B4+C4*C9*C`$11+D4*D9*D`$11)"

$bq.Range("D30").Select() | Out-Null

$focs.Activate() | Out-Null
